$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arrB = New-Object "object[,]" 24,1
$arrB[0,0] = 1.239669154284456
$arrB[1,0] = 1.092238921996511
$arrB[2,0] = 1.001421875712481
$arrB[3,0] = 0.9643416861758851
$arrB[4,0] = 0.9581803041372154
$arrB[5,0] = 1.000922085505977
$arrB[6,0] = 1.188897734072611
$arrB[7,0] = 1.555091194192187
$arrB[8,0] = 1.822559052067106
$arrB[9,0] = 1.943877057524276
$arrB[10,0] = 1.989763988017273
$arrB[11,0] = 1.979883839833974
$arrB[12,0] = 1.947653291195365
$arrB[13,0] = 1.927904097252338
$arrB[14,0] = 1.814623279994976
$arrB[15,0] = 1.745036526342801
$arrB[16,0] = 1.704978855460752
$arrB[17,0] = 1.691410391918907
$arrB[18,0] = 1.752447611123046
$arrB[19,0] = 1.957121658883409
$arrB[20,0] = 2.090574578272708
$arrB[21,0] = 2.019377769559298
$arrB[22,0] = 1.749097219353303
$arrB[23,0] = 1.456295877450543
$ws.Range("B2:B25").Value2 = $arrB

$arrC = New-Object "object[,]" 24,1
$arrC[0,0] = 0.3097235529693023
$arrC[1,0] = 0.2713006013104291
$arrC[2,0] = 0.2476106893304575
$arrC[3,0] = 0.2379329323356103
$arrC[4,0] = 0.2363245248620274
$arrC[5,0] = 0.2474802676092338
$arrC[6,0] = 0.2964960985494542
$arrC[7,0] = 0.3918104739963724
$arrC[8,0] = 0.4613170848028858
$arrC[9,0] = 0.4928185808826697
$arrC[10,0] = 0.5047298738130621
$arrC[11,0] = 0.502165360254196
$arrC[12,0] = 0.4937988879545969
$arrC[13,0] = 0.488671862548415
$arrC[14,0] = 0.4592559604842563
$arrC[15,0] = 0.4411796447118945
$arrC[16,0] = 0.4307716194696809
$arrC[17,0] = 0.4272457737772015
$arrC[18,0] = 0.4431050440267654
$arrC[19,0] = 0.4962568084031886
$arrC[20,0] = 0.5308913950134411
$arrC[21,0] = 0.5124159514994062
$arrC[22,0] = 0.4422346198775244
$arrC[23,0] = 0.3661147615479763
$ws.Range("C2:C25").Value2 = $arrC

$arrD = New-Object "object[,]" 24,1
$arrD[0,0] = 0.4044659625363352
$arrD[1,0] = 0.3916086746908718
$arrD[2,0] = 0.3838530251474168
$arrD[3,0] = 0.3807276203612417
$arrD[4,0] = 0.3802107739315659
$arrD[5,0] = 0.3838107325347551
$arrD[6,0] = 0.4000040747079368
$arrD[7,0] = 0.4328539629383954
$arrD[8,0] = 0.457650548427381
$arrD[9,0] = 0.4690738536491494
$arrD[10,0] = 0.4734200106640571
$arrD[11,0] = 0.4724830838162291
$arrD[12,0] = 0.4694310068832976
$arrD[13,0] = 0.4675641724558375
$arrD[14,0] = 0.4569068737257282
$arrD[15,0] = 0.4504055077931923
$arrD[16,0] = 0.4466795856303918
$arrD[17,0] = 0.4454203755488209
$arrD[18,0] = 0.4510961943832967
$arrD[19,0] = 0.4703269237790266
$arrD[20,0] = 0.4830141299110551
$arrD[21,0] = 0.4762319219245228
$arrD[22,0] = 0.4507838982078454
$arrD[23,0] = 0.4238505802398436
$ws.Range("D2:D25").Value2 = $arrD

$arrF = New-Object "object[,]" 24,1
$arrF[0,0] = 0.8780946014472377
$arrF[1,0] = 0.8813901442152599
$arrF[2,0] = 0.8841514360561291
$arrF[3,0] = 0.8854618446390745
$arrF[4,0] = 0.8856906093244987
$arrF[5,0] = 0.8841683594728238
$arrF[6,0] = 0.8790775468017102
$arrF[7,0] = 0.8749664177868084
$arrF[8,0] = 0.875551645210777
$arrF[9,0] = 0.876606479177056
$arrF[10,0] = 0.8771197473566161
$arrF[11,0] = 0.8770041367203163
$arrF[12,0] = 0.8766464224052299
$arrF[13,0] = 0.8764421480455979
$arrF[14,0] = 0.8754986145034849
$arrF[15,0] = 0.8751220718747845
$arrF[16,0] = 0.8749796970886408
$arrF[17,0] = 0.8749442228429913
$arrF[18,0] = 0.8751544727119551
$arrF[19,0] = 0.8767483992103706
$arrF[20,0] = 0.878453779578507
$arrF[21,0] = 0.8774827216736014
$arrF[22,0] = 0.8751395934928894
$arrF[23,0] = 0.8754472630349426
$ws.Range("F2:F25").Value2 = $arrF

$arrG = New-Object "object[,]" 24,1
$arrG[0,0] = 0.3270373396610253
$arrG[1,0] = 0.3311291648688481
$arrG[2,0] = 0.3340581561467104
$arrG[3,0] = 0.3353560648092824
$arrG[4,0] = 0.3355778695048812
$arrG[5,0] = 0.3340752383836261
$arrG[6,0] = 0.328361475529114
$arrG[7,0] = 0.3204819841344104
$arrG[8,0] = 0.3167473972985562
$arrG[9,0] = 0.3155004086738202
$arrG[10,0] = 0.3150936584514596
$arrG[11,0] = 0.3151783406072397
$arrG[12,0] = 0.3154656305098769
$arrG[13,0] = 0.3156501424235216
$arrG[14,0] = 0.3168380248906857
$arrG[15,0] = 0.3176828321166667
$arrG[16,0] = 0.3182112489553859
$arrG[17,0] = 0.3183974477545704
$arrG[18,0] = 0.3175884982748869
$arrG[19,0] = 0.3153794662289187
$arrG[20,0] = 0.3143174906438873
$arrG[21,0] = 0.3148491998231435
$arrG[22,0] = 0.3176310135516616
$arrG[23,0] = 0.322254777162378
$ws.Range("G2:G25").Value2 = $arrG

$arrH = New-Object "object[,]" 24,1
$arrH[0,0] = 0.5011604010670894
$arrH[1,0] = 0.5078366460444315
$arrH[2,0] = 0.5122855973797087
$arrH[3,0] = 0.5141864097384357
$arrH[4,0] = 0.5145073396181488
$arrH[5,0] = 0.512310876931231
$arrH[6,0] = 0.5033897230123543
$arrH[7,0] = 0.4886749970235442
$arrH[8,0] = 0.4795652223595894
$arrH[9,0] = 0.4757917547059307
$arrH[10,0] = 0.4744162564645649
$arrH[11,0] = 0.4747101160189757
$arrH[12,0] = 0.4756775199890129
$arrH[13,0] = 0.4762770453813019
$arrH[14,0] = 0.4798192968725772
$arrH[15,0] = 0.4820873750854062
$arrH[16,0] = 0.4834267894335582
$arrH[17,0] = 0.4838862778491091
$arrH[18,0] = 0.4818423237760072
$arrH[19,0] = 0.4753919187625257
$arrH[20,0] = 0.4714877130864465
$arrH[21,0] = 0.4735429126582886
$arrH[22,0] = 0.4819530010070423
$arrH[23,0] = 0.4923574590113944
$ws.Range("H2:H25").Value2 = $arrH

$arrJ = New-Object "object[,]" 24,1
$arrJ[0,0] = 0.3465189411931107
$arrJ[1,0] = 0.334809163559072
$arrJ[2,0] = 0.327824504368678
$arrJ[3,0] = 0.3250297399127788
$arrJ[4,0] = 0.3245687827287185
$arrJ[5,0] = 0.327786604583622
$arrJ[6,0] = 0.3424387839639706
$arrJ[7,0] = 0.3728042367256421
$arrJ[8,0] = 0.3961179725474722
$arrJ[9,0] = 0.4069441559018117
$arrJ[10,0] = 0.4110755879714958
$arrJ[11,0] = 0.4101843952293507
$arrJ[12,0] = 0.407283413799874
$arrJ[13,0] = 0.4055106209423798
$arrJ[14,0] = 0.3954149013474506
$arrJ[15,0] = 0.3892780620811038
$arrJ[16,0] = 0.3857690791557218
$arrJ[17,0] = 0.384584562282285
$arrJ[18,0] = 0.3899291895775008
$arrJ[19,0] = 0.4081346390592842
$arrJ[20,0] = 0.4202182862767216
$arrJ[21,0] = 0.413752037662789
$arrJ[22,0] = 0.3896347551392694
$arrJ[23,0] = 0.3644138711846381
$ws.Range("J2:J25").Value2 = $arrJ

$arrO = New-Object "object[,]" 24,1
$arrO[0,0] = 1.588787659320644
$arrO[1,0] = 1.611141893736175
$arrO[2,0] = 1.62647268827287
$arrO[3,0] = 1.633122596567802
$arrO[4,0] = 1.634251083916254
$arrO[5,0] = 1.62656074315494
$arrO[6,0] = 1.596161545051729
$arrO[7,0] = 1.549337562131512
$arrO[8,0] = 1.522804575096515
$arrO[9,0] = 1.512458110483962
$arrO[10,0] = 1.508789272788306
$arrO[11,0] = 1.509568321187061
$arrO[12,0] = 1.512151272336425
$arrO[13,0] = 1.513765888981666
$arrO[14,0] = 1.523515534119994
$arrO[15,0] = 1.529938990353372
$arrO[16,0] = 1.533795741136245
$arrO[17,0] = 1.535129380089273
$arrO[18,0] = 1.52923841123453
$arrO[19,0] = 1.511385824915749
$arrO[20,0] = 1.501170930457477
$arrO[21,0] = 1.506489450989079
$arrO[22,0] = 1.529554632858023
$arrO[23,0] = 1.560628083690133
$ws.Range("O2:O25").Value2 = $arrO

